# NSMB - Begin 8-8
# Adds rows 175-191 (new splits/checkpoints through Enter 8-8 and into 8-8)
# to the "V4" sheet (first worksheet), continuing the existing A/B/C/D table,
# and updates the view state (selection) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row data: Name (column A), B (start frame), C (end frame)
$rowData = @{
    175 = @{ Name = "Checkpoint 791";     B = 56265; C = 66317 }
    176 = @{ Name = "Land on Koopa";      B = 56378; C = 66431 }
    177 = @{ Name = "Checkpoint 1154";    B = 56428; C = 66483 }
    178 = @{ Name = "Checkpoint 1702";    B = 56640; C = 66695 }
    179 = @{ Name = "Checkpoint 1927";    B = 56715; C = 66770 }
    180 = @{ Name = "Checkpoint 2230";    B = 56898; C = 66954 }
    181 = @{ Name = "Checkpoint 2421";    B = 56982; C = 67039 }
    182 = @{ Name = "Checkpoint 2550";    B = 57020; C = 67079 }
    183 = @{ Name = "Checkpoint ";        B = 57085; C = 67144 }
    184 = @{ Name = "Checkpoint 2927";    B = 57134; C = 67193 }
    185 = @{ Name = "Checkpoint 3271";    B = 57234; C = 67293 }
    186 = @{ Name = "Checkpoint 3576";    B = 57309; C = 67368 }
    187 = @{ Name = "Blast out of pipe";  B = 57531; C = 67591 }
    188 = @{ Name = "Get Flag";           B = 57771; C = 67831 }
    189 = @{ Name = "End lLevel";         B = 58289; C = 68349 }
    190 = @{ Name = "Enter 8-8";          B = 58622; C = 69018 }
    191 = @{ Name = "1st Move";           B = 58848; C = 69265 }
}

# Write column A in this exact sequence so that newly-introduced shared
# strings are appended to the shared string table in the same order the
# original author's edit produced (new text is interleaved with rows that
# reuse already-existing labels, and a couple of labels were typed out of
# row order).
$aOrder = @(175, 176, 177, 178, 179, 180, 181, 182, 184, 183, 185, 186, 190, 187, 188, 189, 191)
foreach ($row in $aOrder) {
    $ws.Range("A$row").Value = $rowData[$row].Name
}

# Now fill in the numeric columns and the shared "diff" formula for every
# new row, in row order.
foreach ($row in 175..191) {
    $d = $rowData[$row]
    $ws.Range("B$row").Value = $d.B
    $ws.Range("C$row").Value = $d.C
    $ws.Range("D$row").Formula = "=IF(B$row > 0,C$row-B$row, 0)"
}

# Update selection / frozen-pane view state to the new bottom of the sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 175
$win.ScrollColumn = 1
$ws.Range("B192").Select()
